$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.232.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.480.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.67%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.495.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0990"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "

# Row 11
$ws.Range("E11").Value = "  -0.77%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.98%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.918.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.099.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "

# Row 17
$ws.Range("E17").Value = "  -1.83%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.485.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.09%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "

# Row 22
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.412"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "

# Row 27
$ws.Range("E27").Value = "  -1.19%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "

# Row 35
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "

# Row 37
$ws.Range("E37").Value = "  -0.75%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40
$ws.Range("E40").Value = "  -2.81%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.52%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0911"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.28%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.743.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
